$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "24/10/2025"
$ws.Range("B15").Value = "Csikszereda M. Ciuc"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "Petrolul"
$ws.Range("F15").Value = "D"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1.79
$ws.Range("L15").Value = 2.31
$ws.Range("M15").Value = 11
$ws.Range("N15").Value = 16
$ws.Range("O15").Value = 4
$ws.Range("P15").Value = 3
